# "Generate Report for Handoff"
#
# The localization-status report moved from "In Translation" to
# "Ready for handoff", and the two timestamps that were stamped at
# report-generation time were refreshed a little over a minute later.
# Widening the "Status" columns (they auto-size to fit the new, longer
# status text) follows naturally from the longer label.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" -----------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Refreshed "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" -
$wsOverview.Range("G2").Value = "2016-09-05 10:56:27"
$wsZhCn.Range("H2").Value     = "2016-09-05 10:56:23"
$wsDeDe.Range("H2").Value     = "2016-09-05 10:56:27"

# --- Status columns widen to fit "Ready for handoff" -----------------------
$wsOverview.Columns.Item(5).ColumnWidth = 16.35
$wsOverview.Columns.Item(6).ColumnWidth = 16.35
$wsZhCn.Columns.Item(3).ColumnWidth     = 16.35
$wsDeDe.Columns.Item(3).ColumnWidth     = 16.35
